# Update first three sections of the report:
# split the "Heater/AC-Controller" textbox into two lines:
#   "Human to control" (new first line)
#   "Heater/AC"         (shortened from "Heater/AC-Controller")

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item(12)
$tr = $shape.TextFrame.TextRange

# Insert the new first line ahead of the existing text, using a
# paragraph-mark ("`r") so it becomes its own paragraph. This keeps the
# original run (and its trailing endParaRPr) intact on the 2nd paragraph.
$oldLen = $tr.Length
[void]$tr.InsertBefore("Human to control`r")

# Grab just the old text (now the 2nd paragraph) and shorten it.
$newLen = $tr.Length
$start = $newLen - $oldLen + 1
$old = $tr.Characters($start, $oldLen)
$old.Text = "Heater/AC"
